# Release-Notes.xlsx update
# - Insert a new "Folder Inventory" row at the top (row 2) for a newly
#   detected folder, pushing all existing rows down by one.
# - Refresh the "Metadata" sheet's generation timestamp / folder count /
#   workflow run number.
# - Refresh the "Summary" sheet's folder counts and most-recent-update time.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Folder Inventory sheet: insert new row 2 with the new folder entry
# ---------------------------------------------------------------------
$inventory = $wb.Worksheets.Item("Folder Inventory")

$inventory.Rows.Item(2).Insert()
$inventory.Rows.Item(2).ClearFormats()

$inventory.Range("A2").Value = "Get Started with Microsoft Fabric with Its Lakehouses"
$inventory.Range("B2").Value = "Get Started with Microsoft Fabric with Its Lakehouses"
$inventory.Range("C2").Value = "2025-06-12 16:16:30 +0530"
$inventory.Range("D2").Value = 1
$inventory.Range("E2").Value = "Root"

# ---------------------------------------------------------------------
# 2. Metadata sheet updates
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "2025-06-12 11:08:10 UTC"
$meta.Range("B4").Value = 72

$meta.Range("B5").NumberFormat = "@"
$meta.Range("B5").Value = "11"
$meta.Range("B5").ClearFormats()

# ---------------------------------------------------------------------
# 3. Summary sheet updates
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("B2").Value = 72
$summary.Range("B3").Value = 72
$summary.Range("B5").Value = "2025-06-12 16:16:30 +0530"
